$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9+ down by one
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new call-by-name entry
$ws.Range("A9").Value = "B1"
$ws.Range("B9").Value = "無月"
$ws.Range("C9").Value = "@Shake"

# Leave the final selection on C9, matching the author's last edit location
$ws.Range("C9").Select()
